$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Change the password value from "123456" to "Astegic1!" and clear the explicit style
$ws.Range("B2").Value = "Astegic1!"
$ws.Range("B2").ClearFormats()

# Move the active selection from B3 to B2
$ws.Range("B2").Select()

$wb.Save()
